$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GA")

# Update the revenue value in C2
$ws.Range("C2").Value = 4951

# Move the active cell/selection from D3 to C3
$ws.Activate()
$ws.Range("C3").Select()
